# Apply the updated NATMI TPM results for Spon2-Itga5 ligand-receptor pairs.
# A new "Resolving-Mac" sending/target cluster is introduced (shared string list
# grows from 26->26 unique entries but string usage count increases as rows expand),
# all existing sender/target combination rows get refreshed numeric stats, and four
# new rows are appended for the Resolving-Mac sending cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Spon2"
$ws.Cells.Item(2, 3).Value = "Itga5"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 0.4500866666666667
$ws.Cells.Item(2, 8).Value = 1.35026
$ws.Cells.Item(2, 9).Value = 0.02628438542510526
$ws.Cells.Item(2, 10).Value = 0.02628438542510525
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 35.04689966666667
$ws.Cells.Item(2, 14).Value = 105.140699
$ws.Cells.Item(2, 15).Value = 0.3824629895491901
$ws.Cells.Item(2, 16).Value = 0.3824629895491901
$ws.Cells.Item(2, 17).Value = 15.77414224797111
$ws.Cells.Item(2, 18).Value = 141.96728023174
$ws.Cells.Item(2, 19).Value = 0.01005280462814892
$ws.Cells.Item(2, 20).Value = 0.01005280462814891

# Row 3: ECs -> FAPs
$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Spon2"
$ws.Cells.Item(3, 3).Value = "Itga5"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 0.4500866666666667
$ws.Cells.Item(3, 8).Value = 1.35026
$ws.Cells.Item(3, 9).Value = 0.02628438542510526
$ws.Cells.Item(3, 10).Value = 0.02628438542510525
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 29.913269
$ws.Cells.Item(3, 14).Value = 89.739807
$ws.Cells.Item(3, 15).Value = 0.3264402385872224
$ws.Cells.Item(3, 16).Value = 0.3264402385872223
$ws.Cells.Item(3, 17).Value = 13.46356353331333
$ws.Cells.Item(3, 18).Value = 121.17207179982
$ws.Cells.Item(3, 19).Value = 0.008580281049289871
$ws.Cells.Item(3, 20).Value = 0.008580281049289866

# Row 4: ECs -> MuSCs
$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Spon2"
$ws.Cells.Item(4, 3).Value = "Itga5"
$ws.Cells.Item(4, 4).Value = "MuSCs"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 0.4500866666666667
$ws.Cells.Item(4, 8).Value = 1.35026
$ws.Cells.Item(4, 9).Value = 0.02628438542510526
$ws.Cells.Item(4, 10).Value = 0.02628438542510525
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 8.911727666666666
$ws.Cells.Item(4, 14).Value = 26.735183
$ws.Cells.Item(4, 15).Value = 0.09725271102035077
$ws.Cells.Item(4, 16).Value = 0.09725271102035075
$ws.Cells.Item(4, 17).Value = 4.011049799731111
$ws.Cells.Item(4, 18).Value = 36.09944819758
$ws.Cells.Item(4, 19).Value = 0.002556227740095281
$ws.Cells.Item(4, 20).Value = 0.00255622774009528

# Row 5: ECs -> Resolving-Mac
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Spon2"
$ws.Cells.Item(5, 3).Value = "Itga5"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 0.4500866666666667
$ws.Cells.Item(5, 8).Value = 1.35026
$ws.Cells.Item(5, 9).Value = 0.02628438542510526
$ws.Cells.Item(5, 10).Value = 0.02628438542510525
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 17.76285166666667
$ws.Cells.Item(5, 14).Value = 53.288555
$ws.Cells.Item(5, 15).Value = 0.1938440608432367
$ws.Cells.Item(5, 16).Value = 0.1938440608432367
$ws.Cells.Item(5, 17).Value = 7.994822697144444
$ws.Cells.Item(5, 18).Value = 71.9534042743
$ws.Cells.Item(5, 19).Value = 0.005095072007571188
$ws.Cells.Item(5, 20).Value = 0.005095072007571187

# Row 6: FAPs -> ECs
$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Spon2"
$ws.Cells.Item(6, 3).Value = "Itga5"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 15.76143266666667
$ws.Cells.Item(6, 8).Value = 47.284298
$ws.Cells.Item(6, 9).Value = 0.9204439983318276
$ws.Cells.Item(6, 10).Value = 0.9204439983318274
$ws.Cells.Item(6, 11).Value = 3
$ws.Cells.Item(6, 12).Value = 1
$ws.Cells.Item(6, 13).Value = 35.04689966666667
$ws.Cells.Item(6, 14).Value = 105.140699
$ws.Cells.Item(6, 15).Value = 0.3824629895491901
$ws.Cells.Item(6, 16).Value = 0.3824629895491901
$ws.Cells.Item(6, 17).Value = 552.3893492715891
$ws.Cells.Item(6, 18).Value = 4971.504143444303
$ws.Cells.Item(6, 19).Value = 0.3520357633146005
$ws.Cells.Item(6, 20).Value = 0.3520357633146004

# Row 7: FAPs -> FAPs
$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Spon2"
$ws.Cells.Item(7, 3).Value = "Itga5"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 15.76143266666667
$ws.Cells.Item(7, 8).Value = 47.284298
$ws.Cells.Item(7, 9).Value = 0.9204439983318276
$ws.Cells.Item(7, 10).Value = 0.9204439983318274
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 29.913269
$ws.Cells.Item(7, 14).Value = 89.739807
$ws.Cells.Item(7, 15).Value = 0.3264402385872224
$ws.Cells.Item(7, 16).Value = 0.3264402385872223
$ws.Cells.Item(7, 17).Value = 471.4759751833873
$ws.Cells.Item(7, 18).Value = 4243.283776650486
$ws.Cells.Item(7, 19).Value = 0.3004699584216187
$ws.Cells.Item(7, 20).Value = 0.3004699584216186

# Row 8: FAPs -> MuSCs
$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Spon2"
$ws.Cells.Item(8, 3).Value = "Itga5"
$ws.Cells.Item(8, 4).Value = "MuSCs"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 15.76143266666667
$ws.Cells.Item(8, 8).Value = 47.284298
$ws.Cells.Item(8, 9).Value = 0.9204439983318276
$ws.Cells.Item(8, 10).Value = 0.9204439983318274
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 8.911727666666666
$ws.Cells.Item(8, 14).Value = 26.735183
$ws.Cells.Item(8, 15).Value = 0.09725271102035077
$ws.Cells.Item(8, 16).Value = 0.09725271102035075
$ws.Cells.Item(8, 17).Value = 140.4615955618371
$ws.Cells.Item(8, 18).Value = 1264.154360056534
$ws.Cells.Item(8, 19).Value = 0.08951567418018144
$ws.Cells.Item(8, 20).Value = 0.08951567418018143

# Row 9: FAPs -> Resolving-Mac
$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Spon2"
$ws.Cells.Item(9, 3).Value = "Itga5"
$ws.Cells.Item(9, 4).Value = "Resolving-Mac"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 15.76143266666667
$ws.Cells.Item(9, 8).Value = 47.284298
$ws.Cells.Item(9, 9).Value = 0.9204439983318276
$ws.Cells.Item(9, 10).Value = 0.9204439983318274
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 17.76285166666667
$ws.Cells.Item(9, 14).Value = 53.288555
$ws.Cells.Item(9, 15).Value = 0.1938440608432367
$ws.Cells.Item(9, 16).Value = 0.1938440608432367
$ws.Cells.Item(9, 17).Value = 279.9679905121544
$ws.Cells.Item(9, 18).Value = 2519.71191460939
$ws.Cells.Item(9, 19).Value = 0.1784226024154269
$ws.Cells.Item(9, 20).Value = 0.1784226024154268

# Row 10: MuSCs -> ECs
$ws.Cells.Item(10, 1).Value = "MuSCs"
$ws.Cells.Item(10, 2).Value = "Spon2"
$ws.Cells.Item(10, 3).Value = "Itga5"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 0.9026056666666666
$ws.Cells.Item(10, 8).Value = 2.707817
$ws.Cells.Item(10, 9).Value = 0.05271081546417152
$ws.Cells.Item(10, 10).Value = 0.05271081546417151
$ws.Cells.Item(10, 11).Value = 3
$ws.Cells.Item(10, 12).Value = 1
$ws.Cells.Item(10, 13).Value = 35.04689966666667
$ws.Cells.Item(10, 14).Value = 105.140699
$ws.Cells.Item(10, 15).Value = 0.3824629895491901
$ws.Cells.Item(10, 16).Value = 0.3824629895491901
$ws.Cells.Item(10, 17).Value = 31.63353023823144
$ws.Cells.Item(10, 18).Value = 284.701772144083
$ws.Cells.Item(10, 19).Value = 0.02015993606400272
$ws.Cells.Item(10, 20).Value = 0.02015993606400271

# Row 11: MuSCs -> FAPs
$ws.Cells.Item(11, 1).Value = "MuSCs"
$ws.Cells.Item(11, 2).Value = "Spon2"
$ws.Cells.Item(11, 3).Value = "Itga5"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 0.9026056666666666
$ws.Cells.Item(11, 8).Value = 2.707817
$ws.Cells.Item(11, 9).Value = 0.05271081546417152
$ws.Cells.Item(11, 10).Value = 0.05271081546417151
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 29.913269
$ws.Cells.Item(11, 14).Value = 89.739807
$ws.Cells.Item(11, 15).Value = 0.3264402385872224
$ws.Cells.Item(11, 16).Value = 0.3264402385872223
$ws.Cells.Item(11, 17).Value = 26.99988610792433
$ws.Cells.Item(11, 18).Value = 242.998974971319
$ws.Cells.Item(11, 19).Value = 0.0172069311762512
$ws.Cells.Item(11, 20).Value = 0.0172069311762512

# Row 12: MuSCs -> MuSCs
$ws.Cells.Item(12, 1).Value = "MuSCs"
$ws.Cells.Item(12, 2).Value = "Spon2"
$ws.Cells.Item(12, 3).Value = "Itga5"
$ws.Cells.Item(12, 4).Value = "MuSCs"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 0.9026056666666666
$ws.Cells.Item(12, 8).Value = 2.707817
$ws.Cells.Item(12, 9).Value = 0.05271081546417152
$ws.Cells.Item(12, 10).Value = 0.05271081546417151
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 8.911727666666666
$ws.Cells.Item(12, 14).Value = 26.735183
$ws.Cells.Item(12, 15).Value = 0.09725271102035077
$ws.Cells.Item(12, 16).Value = 0.09725271102035075
$ws.Cells.Item(12, 17).Value = 8.043775891723444
$ws.Cells.Item(12, 18).Value = 72.39398302551099
$ws.Cells.Item(12, 19).Value = 0.00512626970398411
$ws.Cells.Item(12, 20).Value = 0.005126269703984107

# Row 13: MuSCs -> Resolving-Mac
$ws.Cells.Item(13, 1).Value = "MuSCs"
$ws.Cells.Item(13, 2).Value = "Spon2"
$ws.Cells.Item(13, 3).Value = "Itga5"
$ws.Cells.Item(13, 4).Value = "Resolving-Mac"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 0.9026056666666666
$ws.Cells.Item(13, 8).Value = 2.707817
$ws.Cells.Item(13, 9).Value = 0.05271081546417152
$ws.Cells.Item(13, 10).Value = 0.05271081546417151
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 17.76285166666667
$ws.Cells.Item(13, 14).Value = 53.288555
$ws.Cells.Item(13, 15).Value = 0.1938440608432367
$ws.Cells.Item(13, 16).Value = 0.1938440608432367
$ws.Cells.Item(13, 17).Value = 16.03285057049278
$ws.Cells.Item(13, 18).Value = 144.295655134435
$ws.Cells.Item(13, 19).Value = 0.01021767851993349
$ws.Cells.Item(13, 20).Value = 0.01021767851993348

# Row 14: Resolving-Mac -> ECs
$ws.Cells.Item(14, 1).Value = "Resolving-Mac"
$ws.Cells.Item(14, 2).Value = "Spon2"
$ws.Cells.Item(14, 3).Value = "Itga5"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 1
$ws.Cells.Item(14, 6).Value = 0.3333333333333333
$ws.Cells.Item(14, 7).Value = 0.009603
$ws.Cells.Item(14, 8).Value = 0.028809
$ws.Cells.Item(14, 9).Value = 0.0005608007788958107
$ws.Cells.Item(14, 10).Value = 0.0005608007788958106
$ws.Cells.Item(14, 11).Value = 3
$ws.Cells.Item(14, 12).Value = 1
$ws.Cells.Item(14, 13).Value = 35.04689966666667
$ws.Cells.Item(14, 14).Value = 105.140699
$ws.Cells.Item(14, 15).Value = 0.3824629895491901
$ws.Cells.Item(14, 16).Value = 0.3824629895491901
$ws.Cells.Item(14, 17).Value = 0.3365553774990001
$ws.Cells.Item(14, 18).Value = 3.028998397491001
$ws.Cells.Item(14, 19).Value = 0.0002144855424380061
$ws.Cells.Item(14, 20).Value = 0.000214485542438006

# Row 15: Resolving-Mac -> FAPs
$ws.Cells.Item(15, 1).Value = "Resolving-Mac"
$ws.Cells.Item(15, 2).Value = "Spon2"
$ws.Cells.Item(15, 3).Value = "Itga5"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 1
$ws.Cells.Item(15, 6).Value = 0.3333333333333333
$ws.Cells.Item(15, 7).Value = 0.009603
$ws.Cells.Item(15, 8).Value = 0.028809
$ws.Cells.Item(15, 9).Value = 0.0005608007788958107
$ws.Cells.Item(15, 10).Value = 0.0005608007788958106
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 29.913269
$ws.Cells.Item(15, 14).Value = 89.739807
$ws.Cells.Item(15, 15).Value = 0.3264402385872224
$ws.Cells.Item(15, 16).Value = 0.3264402385872223
$ws.Cells.Item(15, 17).Value = 0.287257122207
$ws.Cells.Item(15, 18).Value = 2.585314099863
$ws.Cells.Item(15, 19).Value = 0.0001830679400626486
$ws.Cells.Item(15, 20).Value = 0.0001830679400626485

# Row 16: Resolving-Mac -> MuSCs
$ws.Cells.Item(16, 1).Value = "Resolving-Mac"
$ws.Cells.Item(16, 2).Value = "Spon2"
$ws.Cells.Item(16, 3).Value = "Itga5"
$ws.Cells.Item(16, 4).Value = "MuSCs"
$ws.Cells.Item(16, 5).Value = 1
$ws.Cells.Item(16, 6).Value = 0.3333333333333333
$ws.Cells.Item(16, 7).Value = 0.009603
$ws.Cells.Item(16, 8).Value = 0.028809
$ws.Cells.Item(16, 9).Value = 0.0005608007788958107
$ws.Cells.Item(16, 10).Value = 0.0005608007788958106
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 8.911727666666666
$ws.Cells.Item(16, 14).Value = 26.735183
$ws.Cells.Item(16, 15).Value = 0.09725271102035077
$ws.Cells.Item(16, 16).Value = 0.09725271102035075
$ws.Cells.Item(16, 17).Value = 0.085579320783
$ws.Cells.Item(16, 18).Value = 0.770213887047
$ws.Cells.Item(16, 19).Value = [double]"5.45393960899419e-05"
$ws.Cells.Item(16, 20).Value = [double]"5.453939608994188e-05"

# Row 17: Resolving-Mac -> Resolving-Mac
$ws.Cells.Item(17, 1).Value = "Resolving-Mac"
$ws.Cells.Item(17, 2).Value = "Spon2"
$ws.Cells.Item(17, 3).Value = "Itga5"
$ws.Cells.Item(17, 4).Value = "Resolving-Mac"
$ws.Cells.Item(17, 5).Value = 1
$ws.Cells.Item(17, 6).Value = 0.3333333333333333
$ws.Cells.Item(17, 7).Value = 0.009603
$ws.Cells.Item(17, 8).Value = 0.028809
$ws.Cells.Item(17, 9).Value = 0.0005608007788958107
$ws.Cells.Item(17, 10).Value = 0.0005608007788958106
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 17.76285166666667
$ws.Cells.Item(17, 14).Value = 53.288555
$ws.Cells.Item(17, 15).Value = 0.1938440608432367
$ws.Cells.Item(17, 16).Value = 0.1938440608432367
$ws.Cells.Item(17, 17).Value = 0.170576664555
$ws.Cells.Item(17, 18).Value = 1.535189980995
$ws.Cells.Item(17, 19).Value = 0.0001087079003052141
$ws.Cells.Item(17, 20).Value = 0.000108707900305214
